$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.609.80'
$ws.Range("E2").Value = '  +2.91%  '
$ws.Range("D3").Value = '2.202.03'
$ws.Range("E3").Value = '  +1.92%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.60'
$ws.Range("E5").Value = '  +6.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.613'
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.08'
$ws.Range("E7").Value = '  +3.04%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.587'
$ws.Range("E9").Value = '  +1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.95'
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0919'
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.82'
$ws.Range("E12").Value = '  +2.06%  '
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = '2.534.29'
$ws.Range("E14").Value = '  +2.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.35'
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("D16").Value = '2.195.89'
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.775'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '42.522.96'
$ws.Range("E18").Value = '  +2.77%  '
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.12'
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.93'
$ws.Range("E21").Value = '  +2.85%  '
$ws.Range("E22").Value = '  +9.65%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '228.47'
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.64'
$ws.Range("E24").Value = '  -2.30%  '
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.65'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +2.75%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").Value = '  +1.90%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").Value = '  +2.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.22'
$ws.Range("E30").Value = '  +12.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.01'
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.06'
$ws.Range("E32").Value = '  +1.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0804'
$ws.Range("E33").Value = '  +4.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.17'
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.107'
$ws.Range("E36").Value = '  +2.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.37'
$ws.Range("E37").Value = '  +2.35%  '
$ws.Range("E38").Value = '  +11.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.08'
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.07'
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.26'
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.197'
$ws.Range("E42").Value = '  +4.97%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.23'
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("B44").Value = 'WOONetwork'
$ws.Range("C44").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.483'
$ws.Range("E44").Value = '  +24.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.64'
$ws.Range("E45").Value = '  +7.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.34'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0979'
$ws.Range("E47").Value = '  +2.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.42'
$ws.Range("E48").Value = '  +11.48%  '
$ws.Range("E49").Value = '  +2.82%  '
$ws.Range("E50").Value = '  +2.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.65'
$ws.Range("E51").Value = '  +1.29%  '
